# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.072.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.304.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.980"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.652.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.306.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.152.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "286.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.25%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.03%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.05%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  -7.90%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.48%  "
$ws.Range("E40").Value = "  -6.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +20.83%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.58%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "

Write-Output "Applied all changes"